$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.657.35'
$ws.Range("E2").Value = '  -1.45%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.904.04'
$ws.Range("E3").Value = '  -1.97%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '528.61'
$ws.Range("E5").Value = '  -2.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.66'
$ws.Range("E6").Value = '  -4.73%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.555'
$ws.Range("E8").Value = '  -1.27%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.912.70'
$ws.Range("E9").Value = '  -1.97%  '
$ws.Range("E10").Value = '  -3.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.02'
$ws.Range("E11").Value = '  -1.39%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.364'
$ws.Range("E12").Value = '  -0.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.409.08'
$ws.Range("E13").Value = '  -1.86%  '
$ws.Range("E14").Value = '  +2.43%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.644.09'
$ws.Range("E15").Value = '  -1.59%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.80'
$ws.Range("E16").Value = '  -3.37%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.911.93'
$ws.Range("E17").Value = '  -1.68%  '
$ws.Range("E18").Value = '  -2.69%  '
$ws.Range("E19").Value = '  -1.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.70'
$ws.Range("E20").Value = '  -1.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '362.18'
$ws.Range("E21").Value = '  -4.70%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.65'
$ws.Range("E22").Value = '  +0.26%  '
$ws.Range("E23").Value = '  -0.14%  '
$ws.Range("E24").Value = '  +0.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.81'
$ws.Range("E25").Value = '  -0.24%  '
$ws.Range("E26").Value = '  -2.48%  '
$ws.Range("E27").Value = '  -2.41%  '
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.88'
$ws.Range("E29").Value = '  -4.28%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0860'
$ws.Range("E30").Value = '  -6.93%  '
$ws.Range("E31").Value = '  +0.03%  '
$ws.Range("E32").Value = '  -1.53%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.75'
$ws.Range("E33").Value = '  -2.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '151.71'
$ws.Range("E34").Value = '  -4.32%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.40'
$ws.Range("E35").Value = '  -4.83%  '
$ws.Range("E36").Value = '  -6.73%  '
$ws.Range("E37").Value = '  -4.43%  '
$ws.Range("E38").Value = '  -4.93%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '37.72'
$ws.Range("E39").Value = '  +2.25%  '
$ws.Range("E40").Value = '  -3.18%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.75'
$ws.Range("E41").Value = '  -4.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.292.74'
$ws.Range("E42").Value = '  -4.55%  '
$ws.Range("E43").Value = '  -2.21%  '
$ws.Range("E44").Value = '  -1.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.66'
$ws.Range("E45").Value = '  -6.30%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.998'
$ws.Range("E46").Value = '  +0.13%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.04'
$ws.Range("E47").Value = '  +2.95%  '
$ws.Range("E48").Value = '  -2.52%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.32'
$ws.Range("E49").Value = '  -1.42%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0926'
$ws.Range("E50").Value = '  -2.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '251.55'
$ws.Range("E51").Value = '  -5.05%  '
